# Format the "transferred_at" date column (A) as real dates instead of
# plain text strings, so the values can be used for import/export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the textual "2020-01-01" values in A2:A3 with an actual date
# serial number and apply a yyyy-mm-dd date number format to the cells.
$ws.Range("A2:A3").Value2 = 43831
$ws.Range("A2:A3").NumberFormat = "yyyy\-mm\-dd"

# Move the active selection, matching the author's saved cursor position.
$ws.Range("D25").Select()
